# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "55.890.23"
Set-TextCell "E2" "  -1.86%  "
Set-TextCell "D3" "2.969.56"
Set-TextCell "E3" "  -0.56%  "
Set-TextCell "E4" "  -0.06%  "
Set-TextCell "D5" "502.37"
Set-TextCell "E5" "  +0.43%  "
Set-TextCell "D6" "136.87"
Set-TextCell "E6" "  -0.77%  "
Set-TextCell "E7" "  -0.07%  "
Set-TextCell "E8" "  -1.17%  "
Set-TextCell "D9" "7.14"
Set-TextCell "E9" "  -1.68%  "
Set-TextCell "E10" "  -1.44%  "
Set-TextCell "D11" "0.363"
Set-TextCell "E11" "  +1.20%  "
Set-TextCell "D12" "3.477.13"
Set-TextCell "E13" "  -1.54%  "
Set-TextCell "D14" "25.80"
Set-TextCell "E14" "  -1.33%  "
Set-TextCell "D15" "0.0000160"
Set-TextCell "E15" "  -0.48%  "
Set-TextCell "D16" "55.885.15"
Set-TextCell "E16" "  -2.13%  "
Set-TextCell "B17" "Polkadot"
Set-TextCell "C17" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D17" "5.99"
Set-TextCell "E17" "  -1.54%  "
Set-TextCell "B18" "WrappedEther"
Set-TextCell "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D18" "2.973.58"
Set-TextCell "E18" "  -0.44%  "
Set-TextCell "D19" "12.81"
Set-TextCell "E19" "  +1.37%  "
Set-TextCell "D20" "7.93"
Set-TextCell "E20" "  +0.60%  "
Set-TextCell "D21" "324.98"
Set-TextCell "E21" "  +1.10%  "
Set-TextCell "D22" "1.00"
Set-TextCell "E22" "  +0.04%  "
Set-TextCell "D23" "0.492"
Set-TextCell "E23" "  +0.28%  "
Set-TextCell "D24" "64.35"
Set-TextCell "E24" "  +1.04%  "
Set-TextCell "D25" "3.093.84"
Set-TextCell "E25" "  -0.56%  "
Set-TextCell "E26" "  +0.12%  "
Set-TextCell "D27" "0.163"
Set-TextCell "E27" "  -0.96%  "
Set-TextCell "D28" "0.0₃0910"
Set-TextCell "E28" "  +1.73%  "
Set-TextCell "D29" "6.35"
Set-TextCell "E29" "  -3.03%  "
Set-TextCell "D30" "6.96"
Set-TextCell "E30" "  -1.83%  "
Set-TextCell "D31" "1.77"
Set-TextCell "E31" "  +0.01%  "
Set-TextCell "B32" "Fetch.AI"
Set-TextCell "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D32" "1.15"
Set-TextCell "E32" "  -0.81%  "
Set-TextCell "B33" "EthereumClassic"
Set-TextCell "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D33" "20.12"
Set-TextCell "E33" "  -0.66%  "
Set-TextCell "D34" "152.78"
Set-TextCell "E34" "  -1.62%  "
Set-TextCell "D35" "4.47"
Set-TextCell "E35" "  -2.27%  "
Set-TextCell "D36" "5.73"
Set-TextCell "E36" "  -0.82%  "
Set-TextCell "D37" "25.61"
Set-TextCell "E37" "  +5.36%  "
Set-TextCell "D38" "1.23"
Set-TextCell "E38" "  -1.01%  "
Set-TextCell "D39" "0.0657"
Set-TextCell "E39" "  -1.29%  "
Set-TextCell "D40" "2.999.89"
Set-TextCell "E40" "  -0.63%  "
Set-TextCell "D41" "36.80"
Set-TextCell "E41" "  -2.50%  "
Set-TextCell "E42" "  -0.12%  "
Set-TextCell "D43" "3.77"
Set-TextCell "E43" "  +0.41%  "
Set-TextCell "D44" "0.647"
Set-TextCell "E44" "  +0.40%  "
Set-TextCell "D45" "2.169.27"
Set-TextCell "E45" "  -1.66%  "
Set-TextCell "D46" "1.34"
Set-TextCell "E46" "  -3.23%  "
Set-TextCell "D47" "5.80"
Set-TextCell "E47" "  -2.74%  "
Set-TextCell "D48" "0.917"
Set-TextCell "E48" "  -2.25%  "
Set-TextCell "E49" "  -0.26%  "
Set-TextCell "D50" "19.39"
Set-TextCell "E50" "  +0.34%  "
Set-TextCell "D51" "0.0846"
Set-TextCell "E51" "  -3.49%  "
